$wb = $excel.ActiveWorkbook
$sheets = $wb.Worksheets

# ---------------------------------------------------------------------------
# 1. Update the "总计" (summary) sheet: shift every existing data row down by
#    one and insert the new "2022-Q4" row at the top (row 2). The old last
#    row (2020-Q4) is duplicated one row further down, matching the target.
# ---------------------------------------------------------------------------
$summary = $sheets.Item("总计")

$summaryRows = @(
    @("2022-Q4", 12, 2.32),
    @("2022-Q3", 8, 1.39),
    @("2022-Q2", 16, 3.7),
    @("2022-Q1", 18, 7.32),
    @("2021-Q4", 10, 2.08),
    @("2021-Q3", 12, 2.31),
    @("2021-Q1", 2, 0.03),
    @("2020-Q4", 2, 0.03)
)

# Row 9 is brand new - seed its formatting (column A style) by copying the
# formatting of the row above it before the values are overwritten below.
$summary.Cells.Item(8, 1).Copy($summary.Cells.Item(9, 1))

for ($i = 0; $i -lt $summaryRows.Length; $i++) {
    $r = $i + 2
    $row = $summaryRows[$i]
    $summary.Cells.Item($r, 1).Value = $i
    $summary.Cells.Item($r, 2).Value = $row[0]
    $summary.Cells.Item($r, 3).Value = $row[1]
    $summary.Cells.Item($r, 4).Value = $row[2]
}

# ---------------------------------------------------------------------------
# 2. Create the new "2022-Q4" sheet by copying the structurally-identical
#    "2022-Q3" sheet (same column layout/styles) and placing it before it,
#    then overwrite its contents with the 2022-Q4 fund-holding data.
# ---------------------------------------------------------------------------
$q3 = $sheets.Item("2022-Q3")
$q3.Copy($q3)
$q4 = $sheets.Item("2022-Q3 (2)")
$q4.Name = "2022-Q4"

$q4Rows = @(
    @("012463", "博时成长优势混合A", "16.18", "89.22", "4.38", "0.7087", 7),
    @("012367", "上投摩根安荣回报混合C", "12.39", "25.70", "3.83", "0.4745", 1),
    @("012366", "上投摩根安荣回报混合A", "9.72", "25.70", "3.83", "0.3723", 1),
    @("004823", "上投摩根安裕回报混合A", "4.26", "36.12", "3.91", "0.1666", 1),
    @("004824", "上投摩根安裕回报混合C", "3.64", "36.12", "3.91", "0.1423", 1),
    @("011034", "南方宝恒混合C", "12.75", "28.78", "1.00", "0.1275", 5),
    @("011033", "南方宝恒混合A", "12.74", "28.78", "1.00", "0.1274", 5),
    @("010742", "南方宁悦一年持有期混合A", "11.05", "28.20", "0.80", "0.0884", 8),
    @("016174", "汇丰晋信策略优选混合A", "1.78", "74.92", "2.44", "0.0434", 6),
    @("012464", "博时成长优势混合C", "0.70", "89.22", "4.38", "0.0307", 7),
    @("010743", "南方宁悦一年持有期混合C", "2.93", "28.20", "0.80", "0.0234", 8),
    @("016175", "汇丰晋信策略优选混合C", "0.40", "74.92", "2.44", "0.0098", 6)
)

# Rows 10-13 are new (the copied sheet only had 9 rows); seed column A's
# formatting for them by copying down from row 9 before values are set.
for ($r = 10; $r -le 13; $r++) {
    $q4.Cells.Item(9, 1).Copy($q4.Cells.Item($r, 1))
}

for ($i = 0; $i -lt $q4Rows.Length; $i++) {
    $r = $i + 2
    $row = $q4Rows[$i]
    $q4.Cells.Item($r, 1).Value = $i
    $q4.Cells.Item($r, 2).Value = $row[0]
    $q4.Cells.Item($r, 3).Value = $row[1]
    $q4.Cells.Item($r, 4).Value = $row[2]
    $q4.Cells.Item($r, 5).Value = $row[3]
    $q4.Cells.Item($r, 6).Value = $row[4]
    $q4.Cells.Item($r, 7).Value = $row[5]
    $q4.Cells.Item($r, 8).Value = $row[6]
}
